$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45047
$ws.Range("B2").Value = 1298
$ws.Range("C2").Value = 12
